$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("D2", "E2", "G2", "D3", "E3", "G3", "D4", "E4", "G4", "D5", "E5", "G5", "D6", "E6", "G6", "D7", "E7", "G7", "D8", "E8", "G8", "D9", "E9", "G9", "D10", "E10", "G10", "D11", "E11", "G11", "E12", "G12", "D13", "E13", "G13", "D14", "E14", "G14", "D15", "E15", "G15", "D16", "E16", "G16", "D17", "E17", "G17", "E18", "G18", "D19", "E19", "G19", "D20", "E20", "G20", "D21", "E21", "G21", "D22", "E22", "G22", "D23", "E23", "G23", "E24", "G24", "D25", "E25", "G25", "D26", "E26", "G26", "E27", "G27", "E28", "G28", "G29", "G30", "G31", "G32", "G33", "G34", "G35", "G36", "G37", "G38", "G39", "D40", "E40", "G40", "D41", "E41", "G41", "D42", "E42", "G42", "D43", "E43", "G43", "D44", "E44", "G44", "E45", "G45", "E46", "G46", "D47", "E47", "G47", "D48", "E48", "G48", "D49", "E49", "G49", "D50", "E50", "G50", "G51")
foreach ($addr in $cells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "255.38"
$ws.Range("E2").Value = "3.93%"
$ws.Range("G2").Value = "14"
$ws.Range("D3").Value = "28.17"
$ws.Range("E3").Value = "-5.36%"
$ws.Range("G3").Value = "14"
$ws.Range("D4").Value = "5.198"
$ws.Range("E4").Value = "-2.21%"
$ws.Range("G4").Value = "14"
$ws.Range("D5").Value = "0.05857"
$ws.Range("E5").Value = "1.86%"
$ws.Range("G5").Value = "14"
$ws.Range("D6").Value = "6.726"
$ws.Range("E6").Value = "1.13%"
$ws.Range("G6").Value = "14"
$ws.Range("D7").Value = "0.8699"
$ws.Range("E7").Value = "1.33%"
$ws.Range("G7").Value = "14"
$ws.Range("D8").Value = "0.9547"
$ws.Range("E8").Value = "11.40%"
$ws.Range("G8").Value = "14"
$ws.Range("D9").Value = "0.1410"
$ws.Range("E9").Value = "2.05%"
$ws.Range("G9").Value = "14"
$ws.Range("D10").Value = "0.07165"
$ws.Range("E10").Value = "1.11%"
$ws.Range("G10").Value = "14"
$ws.Range("D11").Value = "0.03211"
$ws.Range("E11").Value = "2.08%"
$ws.Range("G11").Value = "14"
$ws.Range("E12").Value = "-1.33%"
$ws.Range("G12").Value = "14"
$ws.Range("D13").Value = "0.001538"
$ws.Range("E13").Value = "0.87%"
$ws.Range("G13").Value = "14"
$ws.Range("D14").Value = "0.0006101"
$ws.Range("E14").Value = "2.28%"
$ws.Range("G14").Value = "14"
$ws.Range("D15").Value = "0.005836"
$ws.Range("E15").Value = "-2.75%"
$ws.Range("G15").Value = "14"
$ws.Range("D16").Value = "3.498"
$ws.Range("E16").Value = "-0.53%"
$ws.Range("G16").Value = "14"
$ws.Range("D17").Value = "3.235"
$ws.Range("E17").Value = "-0.26%"
$ws.Range("G17").Value = "14"
$ws.Range("E18").Value = "1.48%"
$ws.Range("G18").Value = "14"
$ws.Range("D19").Value = "0.3177"
$ws.Range("E19").Value = "0.69%"
$ws.Range("G19").Value = "14"
$ws.Range("D20").Value = "0.03453"
$ws.Range("E20").Value = "4.13%"
$ws.Range("G20").Value = "14"
$ws.Range("D21").Value = "0.1309"
$ws.Range("E21").Value = "0.83%"
$ws.Range("G21").Value = "14"
$ws.Range("D22").Value = "3.532"
$ws.Range("E22").Value = "1.61%"
$ws.Range("G22").Value = "14"
$ws.Range("D23").Value = "0.04186"
$ws.Range("E23").Value = "1.94%"
$ws.Range("G23").Value = "14"
$ws.Range("E24").Value = "-0.71%"
$ws.Range("G24").Value = "14"
$ws.Range("D25").Value = "0.001227"
$ws.Range("E25").Value = "0.34%"
$ws.Range("G25").Value = "14"
$ws.Range("D26").Value = "0.004568"
$ws.Range("E26").Value = "9.51%"
$ws.Range("G26").Value = "14"
$ws.Range("E27").Value = "0.11%"
$ws.Range("G27").Value = "14"
$ws.Range("E28").Value = "1.29%"
$ws.Range("G28").Value = "14"
$ws.Range("G29").Value = "14"
$ws.Range("G30").Value = "14"
$ws.Range("G31").Value = "14"
$ws.Range("G32").Value = "14"
$ws.Range("G33").Value = "14"
$ws.Range("G34").Value = "14"
$ws.Range("G35").Value = "14"
$ws.Range("G36").Value = "14"
$ws.Range("G37").Value = "14"
$ws.Range("G38").Value = "14"
$ws.Range("G39").Value = "14"
$ws.Range("D40").Value = "0.03823"
$ws.Range("E40").Value = "1.34%"
$ws.Range("G40").Value = "14"
$ws.Range("D41").Value = "0.005631"
$ws.Range("E41").Value = "56.49%"
$ws.Range("G41").Value = "14"
$ws.Range("D42").Value = "0.1104"
$ws.Range("E42").Value = "2.91%"
$ws.Range("G42").Value = "14"
$ws.Range("D43").Value = "0.002344"
$ws.Range("E43").Value = "-2.02%"
$ws.Range("G43").Value = "14"
$ws.Range("D44").Value = "0.009834"
$ws.Range("E44").Value = "7.07%"
$ws.Range("G44").Value = "14"
$ws.Range("E45").Value = "1.49%"
$ws.Range("G45").Value = "14"
$ws.Range("E46").Value = "0.09%"
$ws.Range("G46").Value = "14"
$ws.Range("D47").Value = "0.09001"
$ws.Range("E47").Value = "0.22%"
$ws.Range("G47").Value = "14"
$ws.Range("D48").Value = "0.002127"
$ws.Range("E48").Value = "-26.63%"
$ws.Range("G48").Value = "14"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").Value = "0.09%"
$ws.Range("G49").Value = "14"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").Value = "0.09%"
$ws.Range("G50").Value = "14"
$ws.Range("G51").Value = "14"
